# Restore C10 on the "Rules" sheet from 18 to 1, per the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
